$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37625
$ws.Range("D2").Value = 54420936
$ws.Range("C3").Value = 90757
$ws.Range("D3").Value = 133049076
$ws.Range("C4").Value = 31118
$ws.Range("D4").Value = 46085954
$ws.Range("C5").Value = 8677
$ws.Range("D5").Value = 12897203
$ws.Range("C6").Value = 1985
$ws.Range("D6").Value = 2950006
$ws.Range("C11").Value = 41195
$ws.Range("D11").Value = 55906501
$ws.Range("C12").Value = 9630
$ws.Range("D12").Value = 13928670
$ws.Range("C13").Value = 25894
$ws.Range("D13").Value = 37977475
$ws.Range("C14").Value = 8302
$ws.Range("D14").Value = 12321218
$ws.Range("C15").Value = 2147
$ws.Range("D15").Value = 3192883
$ws.Range("C16").Value = 416
$ws.Range("D16").Value = 613123
$ws.Range("C19").Value = 10196
$ws.Range("D19").Value = 13504149
$ws.Range("C20").Value = 13353
$ws.Range("D20").Value = 19283099
$ws.Range("C21").Value = 31593
$ws.Range("D21").Value = 46366366
$ws.Range("C23").Value = 2633
$ws.Range("D23").Value = 3914682
$ws.Range("C26").Value = 11656
$ws.Range("D26").Value = 15571604
$ws.Range("C27").Value = 7625
$ws.Range("D27").Value = 11047098
$ws.Range("C28").Value = 22437
$ws.Range("D28").Value = 32934253
$ws.Range("C29").Value = 7795
$ws.Range("D29").Value = 11600133
$ws.Range("C30").Value = 1956
$ws.Range("D30").Value = 2918499
$ws.Range("C31").Value = 367
$ws.Range("D31").Value = 547915
$ws.Range("C33").Value = 8276
$ws.Range("D33").Value = 10937183
$ws.Range("C34").Value = 3234
$ws.Range("D34").Value = 4668137
$ws.Range("C35").Value = 7806
$ws.Range("D35").Value = 11399456
$ws.Range("C36").Value = 3170
$ws.Range("D36").Value = 4697961
$ws.Range("C37").Value = 827
$ws.Range("D37").Value = 1231723
$ws.Range("C40").Value = 2460
$ws.Range("D40").Value = 3324831
$ws.Range("C41").Value = 17189
$ws.Range("D41").Value = 24857206
$ws.Range("C42").Value = 50973
$ws.Range("D42").Value = 74729656
$ws.Range("C43").Value = 18973
$ws.Range("D43").Value = 28181943
$ws.Range("C44").Value = 5594
$ws.Range("D44").Value = 8330978
$ws.Range("C45").Value = 1199
$ws.Range("D45").Value = 1789045
$ws.Range("C49").Value = 16636
$ws.Range("D49").Value = 22157416
$ws.Range("C50").Value = 2000
$ws.Range("D50").Value = 2901262
$ws.Range("C51").Value = 6842
$ws.Range("D51").Value = 10057924
$ws.Range("C53").Value = 751
$ws.Range("D53").Value = 1121805
$ws.Range("C56").Value = 6832
$ws.Range("D56").Value = 9406488
$ws.Range("C57").Value = 931
$ws.Range("D57").Value = 1366079
$ws.Range("C58").Value = 2340
$ws.Range("D58").Value = 3468917
$ws.Range("C59").Value = 935
$ws.Range("D59").Value = 1392001
$ws.Range("C63").Value = 1372
$ws.Range("D63").Value = 1929206
$ws.Range("C64").Value = 15303
$ws.Range("D64").Value = 22107537
$ws.Range("C65").Value = 44580
$ws.Range("D65").Value = 65238740
$ws.Range("C66").Value = 15672
$ws.Range("D66").Value = 23291182
$ws.Range("C67").Value = 4561
$ws.Range("D67").Value = 6793792
$ws.Range("C68").Value = 917
$ws.Range("D68").Value = 1363668
$ws.Range("C72").Value = 15046
$ws.Range("D72").Value = 19844168
$ws.Range("C73").Value = 51141
$ws.Range("D73").Value = 74423738
$ws.Range("C74").Value = 145431
$ws.Range("D74").Value = 214259690
$ws.Range("C75").Value = 63411
$ws.Range("D75").Value = 94492270
$ws.Range("C76").Value = 20251
$ws.Range("D76").Value = 30257186
$ws.Range("C77").Value = 4790
$ws.Range("D77").Value = 7156723
$ws.Range("C84").Value = 50585
$ws.Range("D84").Value = 68827969
$ws.Range("C85").Value = 4569
$ws.Range("D85").Value = 6619511
$ws.Range("C86").Value = 11528
$ws.Range("D86").Value = 16936635
$ws.Range("C87").Value = 3869
$ws.Range("D87").Value = 5766083
$ws.Range("C92").Value = 5375
$ws.Range("D92").Value = 7229085
$ws.Range("C93").Value = 1588
$ws.Range("D93").Value = 2286932
$ws.Range("C94").Value = 5129
$ws.Range("D94").Value = 7554213
$ws.Range("C100").Value = 3524
$ws.Range("D100").Value = 4664149
$ws.Range("C101").Value = 596
$ws.Range("D101").Value = 887664
$ws.Range("C106").Value = 10725
$ws.Range("D106").Value = 15559883
$ws.Range("C107").Value = 29117
$ws.Range("D107").Value = 42783513
$ws.Range("C108").Value = 9758
$ws.Range("D108").Value = 14510650
$ws.Range("C113").Value = 9762
$ws.Range("D113").Value = 12898003
$ws.Range("C114").Value = 30321
$ws.Range("D114").Value = 43727252
$ws.Range("C115").Value = 65974
$ws.Range("D115").Value = 96559730
$ws.Range("C116").Value = 21312
$ws.Range("D116").Value = 31672832
$ws.Range("C117").Value = 6046
$ws.Range("D117").Value = 9007561
$ws.Range("C123").Value = 25761
$ws.Range("D123").Value = 34421398
$ws.Range("C124").Value = 35846
$ws.Range("D124").Value = 51743035
$ws.Range("C125").Value = 76557
$ws.Range("D125").Value = 111957399
$ws.Range("C126").Value = 23778
$ws.Range("D126").Value = 35292600
$ws.Range("C127").Value = 6373
$ws.Range("D127").Value = 9470551
$ws.Range("C132").Value = 31665
$ws.Range("D132").Value = 42065299
$ws.Range("C133").Value = 13179
$ws.Range("D133").Value = 19077712
$ws.Range("C134").Value = 32227
$ws.Range("D134").Value = 47337706
$ws.Range("C135").Value = 11453
$ws.Range("D135").Value = 17017542
$ws.Range("C136").Value = 2950
$ws.Range("D136").Value = 4397805
$ws.Range("C137").Value = 495
$ws.Range("D137").Value = 736490
$ws.Range("C140").Value = 10782
$ws.Range("D140").Value = 14380641
$ws.Range("C141").Value = 34907
$ws.Range("D141").Value = 50413322
$ws.Range("C142").Value = 80995
$ws.Range("D142").Value = 118675781
$ws.Range("C143").Value = 24283
$ws.Range("D143").Value = 36079055
$ws.Range("C144").Value = 6367
$ws.Range("D144").Value = 9500067
$ws.Range("C145").Value = 1429
$ws.Range("D145").Value = 2125730
$ws.Range("C147").Value = 14
$ws.Range("D147").Value = 21000
$ws.Range("C148").Value = 29071
$ws.Range("D148").Value = 39231650
